# Add the new worksheet '09081' right after '09810' (the last existing sheet)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "09081"

# Whole used range (A1:G67) is formatted as Text so numeric-looking values
# (e.g. '10142', '80079') are preserved as text, matching the source form.
$ws.Range("A1:G67").NumberFormat = "@"

# Row data: omschrijving | inhoud | weergave | uitlijnen | regel verwijderen | regel template | A327
$data = @(
  ,@('omschrijving', 'inhoud', 'weergave', 'uitlijnen', 'regel verwijderen', 'regel template', 'A327')
  ,@('Verzekerde Modules', $null, $null, $null, $null, ([char]0x0c + '35 Verzekerde Modules'), 'x')
  ,@('Soort verzekering', '10142', $null, $null, $null, '02 Soort verzekering         10142', 'x')
  ,@($null, '10142', 'Omschrijving', 'Links', 'niet verwijderen', $null, 'x')
  ,@('Gezinssamenstelling', '10694', $null, $null, $null, '03 Gezinssamenstelling       10694', 'x')
  ,@($null, '10694', 'Omschrijving', 'Links', 'verwijderen', $null, 'x')
  ,@('Aantal artsen', '13610', $null, $null, $null, '04 Aantal artsen             13610', 'x')
  ,@($null, '13610', 'Getal exclusief decimalen', 'Links', 'verwijderen', $null, 'x')
  ,@('n/a', 'Vrije advocaat keuze (VAK)', $null, $null, $null, '06                           Vrije advocaat keuze (VAK)', 'x')
  ,@('Verzekerd bedrag VAK', '€ 10611', $null, $null, $null, '07 Verzekerd bedrag VAK      € 10611', 'x')
  ,@($null, '10611', 'Getal inclusief decimalen', 'Links', 'verwijderen', $null, 'x')
  ,@('Eigen bijdrage VAK', '€ 13616', $null, $null, $null, '08 Eigen bijdrage VAK        € 13616', 'x')
  ,@($null, '13616', 'Getal inclusief decimalen', 'Links', 'verwijderen', $null, 'x')
  ,@('Militair', 'Ja                                                      81005', $null, $null, $null, '09 Militair                  Ja                                                      81005', 'x')
  ,@($null, '81005', 'Getal exclusief decimalen', 'Links', 'verwijderen', $null, 'x')
  ,@('Lid van militaire vakbond', '80991                                                       81005', $null, $null, $null, '10 Lid van militaire vakbond 80991                                                       81005', 'x')
  ,@($null, '80991', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '81005', 'Getal exclusief decimalen', 'Links', 'verwijderen', $null, 'x')
  ,@('Onroerend goed object 1', '80079 80095 86578', $null, $null, $null, '12 Onroerend goed object 1   80079 80095 86578', 'x')
  ,@($null, '80079', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80095', 'Getal exclusief decimalen', 'Links', 'niet verwijderen', $null, 'x')
  ,@($null, '86578', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('n/a', '80084 80090', $null, $null, $null, '13                           80084 80090', 'x')
  ,@($null, '80084', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80090', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('Dekkingscombinatie', '80732', $null, $null, $null, '14 Dekkingscombinatie        80732', 'x')
  ,@($null, '80732', 'Omschrijving', 'Links', 'verwijderen', $null, 'x')
  ,@('Herbouwwaarde', '€ 80008', $null, $null, $null, '15 Herbouwwaarde             € 80008', 'x')
  ,@($null, '80008', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('Huurwaarde', '€ 84818', $null, $null, $null, '16 Huurwaarde                € 84818', 'x')
  ,@($null, '84818', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('n/a', '80079', $null, $null, $null, '17                                                                                   80079', 'x')
  ,@($null, '80079', $null, 'Links', 'verwijderen', $null, 'x')
  ,@('Onroerend goed object 2', '80077 80093 86451', $null, $null, $null, '18 Onroerend goed object 2   80077 80093 86451', 'x')
  ,@($null, '80077', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80093', 'Getal exclusief decimalen', 'Links', 'niet verwijderen', $null, 'x')
  ,@($null, '86451', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('n/a', '80082 80088', $null, $null, $null, '19                           80082 80088', 'x')
  ,@($null, '80082', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80088', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('Dekkingscombinatie', '80733', $null, $null, $null, '20 Dekkingscombinatie        80733', 'x')
  ,@($null, '80733', 'Omschrijving', 'Links', 'verwijderen', $null, 'x')
  ,@('Herbouwwaarde', '€ 80735', $null, $null, $null, '21 Herbouwwaarde             € 80735', 'x')
  ,@($null, '80735', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('Huurwaarde', '€ 80737', $null, $null, $null, '22 Huurwaarde                € 80737', 'x')
  ,@($null, '80737', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('n/a', '80077', $null, $null, $null, '23                                                                                   80077', 'x')
  ,@($null, '80077', $null, 'Links', 'verwijderen', $null, 'x')
  ,@('Onroerend goed object 3', '80078 80094 86577', $null, $null, $null, '24 Onroerend goed object 3   80078 80094 86577', 'x')
  ,@($null, '80078', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80094', 'Getal exclusief decimalen', 'Links', 'niet verwijderen', $null, 'x')
  ,@($null, '86577', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('n/a', '80083 80089', $null, $null, $null, '25                           80083 80089', 'x')
  ,@($null, '80083', $null, 'Links', 'verwijderen', $null, 'x')
  ,@($null, '80089', $null, 'Links', 'niet verwijderen', $null, 'x')
  ,@('Dekkingscombinatie', '80734', $null, $null, $null, '26 Dekkingscombinatie        80734', 'x')
  ,@($null, '80734', 'Omschrijving', 'Links', 'verwijderen', $null, 'x')
  ,@('Herbouwwaarde', '€ 80736', $null, $null, $null, '27 Herbouwwaarde             € 80736', 'x')
  ,@($null, '80736', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('Huurwaarde', '€ 80738', $null, $null, $null, '28 Huurwaarde                € 80738', 'x')
  ,@($null, '80738', 'Getal inclusief decimalen', 'Rechts', 'verwijderen', $null, 'x')
  ,@('n/a', '80078', $null, $null, $null, '29                                                                                   80078', 'x')
  ,@($null, '80078', $null, 'Links', 'verwijderen', $null, 'x')
  ,@('Geschil', 'Als u een juridisch geschil wilt melden of behoefte heeft aan juridisch advies kunt u', $null, $null, $null, '30 Geschil                   Als u een juridisch geschil wilt melden of behoefte heeft aan juridisch advies kunt u', 'x')
  ,@('n/a', 'contact opnemen met:', $null, $null, $null, '31                           contact opnemen met:', 'x')
  ,@('n/a', 'ARAG Rechtsbijstand', $null, $null, $null, '32                           ARAG Rechtsbijstand', 'x')
  ,@('n/a', 'T (033) - 434 23 42 of via www.ARAG.nl', $null, $null, $null, '33                           T (033) - 434 23 42 of via www.ARAG.nl', 'x')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 1
    for ($j = 0; $j -lt $row.Count; $j++) {
        $val = $row[$j]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $j + 1).Value = $val
        }
    }
}

